$wb = $excel.ActiveWorkbook

# Delete the "Desarquivamentos Pendentes" sheet entirely
$wb.Worksheets("Desarquivamentos Pendentes").Delete()

# Rename remaining sheets to their updated (uppercase) names
$wb.Worksheets("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Keep the original active sheet selected (first tab), since deleting a
# later sheet would otherwise shift the active tab
$wb.Worksheets("PAINEIS DARQ").Select()
